$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-05 Monday", "2026-01-06 Tuesday"),
    @("34×35=1190", "62×91=5642"),
    @("28×15=420", "98×53=5194"),
    @("84×97=8148", "83×62=5146"),
    @("76×25=1900", "35×73=2555"),
    @("11×25=275", "42×44=1848"),
    @("40×56=2240", "41×12=492"),
    @("87×50=4350", "26×19=494"),
    @("12×67=804", "91×63=5733"),
    @("11×21=231", "40×96=3840"),
    @("81×62=5022", "21×96=2016"),
    @("79×81=6399", "50×25=1250"),
    @("51×96=4896", "46×36=1656"),
    @("19×68=1292", "83×19=1577"),
    @("13×13=169", "23×62=1426"),
    @("15×61=915", "61×87=5307"),
    @("77×27=2079", "89×60=5340"),
    @("36×82=2952", "51×92=4692"),
    @("76×42=3192", "18×99=1782"),
    @("82×88=7216", "77×99=7623"),
    @("25×75=1875", "85×80=6800"),
    @("24×58=1392", "40×76=3040"),
    @("24×56=1344", "44×65=2860"),
    @("27×12=324", "16×27=432"),
    @("25×74=1850", "65×88=5720"),
    @("33×25=825", "82×18=1476")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
